$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.299.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.51%  "

$ws.Range("E4").Value = "  -0.46%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4464"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3816"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07835"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.009"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.843.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.838"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.077"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.86%  "

$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.99%  "

$ws.Range("E18").Value = "  -3.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06504"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.60%  "

$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.451"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.31%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.295.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.248"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.056.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "

$ws.Range("E28").Value = "  -4.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.050"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.509"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09297"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.454"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9253"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.617"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.228"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.52%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02206"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05937"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.202"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.269"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5882"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1846"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.91%  "

$ws.Range("E44").Value = "  -8.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.249"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5620"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.353"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.913"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06849"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.90%  "
